# Add an "Apêndices" section (with Modelos / Glossários / Índice sub-items
# and a trailing blank line) right after the existing "Requisitos de
# qualidade" bullet, before the document's final empty paragraph.

$d = $word.ActiveDocument

# Locate the anchor paragraph ("Requisitos de qualidade") robustly via Find
# rather than a hard-coded paragraph index.
$rng = $d.Content
$found = $rng.Find.Execute("Requisitos de qualidade", $true, $false, $false,
                            $false, $false, $true, 1, $false, "", 0)
$anchorIndex = $rng.Paragraphs.First.Index

# The document always ends with a trailing empty paragraph; insert the new
# content immediately before it so it lands right after the anchor above.
$tailRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$tailRange.Collapse(1)
$tailRange.InsertBefore("Apêndices`rModelos`rGlossários`rÍndice`r")

# Add the extra blank paragraph that follows "Índice" in the target layout.
$d.Paragraphs.Item($d.Paragraphs.Count).Range.InsertParagraphBefore()

# The five freshly-created paragraphs now sit right after the anchor.
$pApendices  = $d.Paragraphs.Item($anchorIndex + 1)
$pModelos    = $d.Paragraphs.Item($anchorIndex + 2)
$pGlossarios = $d.Paragraphs.Item($anchorIndex + 3)
$pIndice     = $d.Paragraphs.Item($anchorIndex + 4)
$pBlank      = $d.Paragraphs.Item($anchorIndex + 5)

# Match the target indentation: 708 twips == 35.4 points for the two
# "Modelos"/"Glossários" sub-items, and an explicit zero first-line indent
# for "Índice" and the trailing blank paragraph.
$pModelos.Format.FirstLineIndent    = 35.4
$pGlossarios.Format.FirstLineIndent = 35.4
$pIndice.Format.FirstLineIndent     = 0
$pBlank.Format.FirstLineIndent      = 0

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
